$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 1087.2142
$ws.Range("J49").Value = 1182.3334
$ws.Range("L49").Value = 3547.0002
$ws.Range("N49").Value = -3819.0002
$ws.Range("H61").Value = 324.4
$ws.Range("I61").Value = 176.25
$ws.Range("K61").Value = 528.75
$ws.Range("M61").Value = -356.75
$ws.Range("H74").Value = 3980.7083
$ws.Range("I74").Value = 4301.7
$ws.Range("J74").Value = 3751.4285
$ws.Range("K74").Value = 4301.7
$ws.Range("L74").Value = 3751.4285
$ws.Range("M74").Value = -3365.7
$ws.Range("N74").Value = -5623.4285
$ws.Range("H77").Value = 3980.7083
$ws.Range("I77").Value = 4301.7
$ws.Range("J77").Value = 3751.4285
$ws.Range("K77").Value = 21508.5
$ws.Range("L77").Value = 18757.1425
$ws.Range("M77").Value = -16828.5
$ws.Range("N77").Value = -28117.1425

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6923.8
$ws.Range("I32").Value = 2846.3896
$ws.Range("J32").Value = 20574.262
$ws.Range("K32").Value = 2846.3896
$ws.Range("L32").Value = 20574.262
$ws.Range("M32").Value = -2559.3896
$ws.Range("N32").Value = -21148.262
$ws.Range("H61").Value = 1208.4667
$ws.Range("I61").Value = 1098.8235
$ws.Range("J61").Value = 1351.8462
$ws.Range("K61").Value = 1098.8235
$ws.Range("L61").Value = 1351.8462
$ws.Range("M61").Value = -886.8235
$ws.Range("N61").Value = -1775.8462
$ws.Range("H88").Value = 90911410
$ws.Range("I88").Value = 999.5
$ws.Range("J88").Value = 111113720
$ws.Range("K88").Value = 999.5
$ws.Range("L88").Value = 111113720
$ws.Range("M88").Value = -593.5
$ws.Range("N88").Value = -111114532
$ws.Range("H91").Value = 90911410
$ws.Range("I91").Value = 999.5
$ws.Range("J91").Value = 111113720
$ws.Range("K91").Value = 999.5
$ws.Range("L91").Value = 111113720
$ws.Range("M91").Value = 404.5
$ws.Range("N91").Value = -111116528
$ws.Range("H104").Value = 21919.334
$ws.Range("J104").Value = 21919.334
$ws.Range("L104").Value = 21919.334
$ws.Range("N104").Value = -28907.334
$ws.Range("H132").Value = 1601.66
$ws.Range("I132").Value = 1293.8857
$ws.Range("J132").Value = 2319.8
$ws.Range("K132").Value = 3881.6571
$ws.Range("L132").Value = 6959.400000000001
$ws.Range("M132").Value = -1351.6571
$ws.Range("N132").Value = -12019.4
$ws.Range("H136").Value = 1208.4667
$ws.Range("I136").Value = 1098.8235
$ws.Range("J136").Value = 1351.8462
$ws.Range("K136").Value = 3296.4705
$ws.Range("L136").Value = 4055.5386
$ws.Range("M136").Value = -746.4704999999999
$ws.Range("N136").Value = -9155.5386

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1403.1538
$ws.Range("I99").Value = 1640
$ws.Range("J99").Value = 1255.125
$ws.Range("K99").Value = 1640
$ws.Range("L99").Value = 1255.125
$ws.Range("M99").Value = -142
$ws.Range("N99").Value = -4251.125
$ws.Range("H134").Value = 1363.6316
$ws.Range("I134").Value = 1397.9166
$ws.Range("J134").Value = 1304.8572
$ws.Range("K134").Value = 4193.7498
$ws.Range("L134").Value = 3914.5716
$ws.Range("M134").Value = -1658.7498
$ws.Range("N134").Value = -8984.571599999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 24714.334
$ws.Range("J28").Value = 24714.334
$ws.Range("L28").Value = 24714.334
$ws.Range("N28").Value = -25204.334
$ws.Range("H96").Value = 19411.154
$ws.Range("J96").Value = 19411.154
$ws.Range("L96").Value = 19411.154
$ws.Range("N96").Value = -24903.154
$ws.Range("H134").Value = 3685.4348
$ws.Range("I134").Value = 5033.615
$ws.Range("J134").Value = 1932.8
$ws.Range("K134").Value = 15100.845
$ws.Range("L134").Value = 5798.4
$ws.Range("M134").Value = -12565.845
$ws.Range("N134").Value = -10868.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 499.04443
$ws.Range("I5").Value = 426
$ws.Range("J5").Value = 836.875
$ws.Range("K5").Value = 1278
$ws.Range("L5").Value = 2510.625
$ws.Range("M5").Value = -1166
$ws.Range("N5").Value = -2734.625
$ws.Range("H98").Value = 521.1667
$ws.Range("I98").Value = 494.44446
$ws.Range("K98").Value = 1483.33338
$ws.Range("M98").Value = 14.66661999999997
$ws.Range("H107").Value = 204.12
$ws.Range("J107").Value = 232.18182
$ws.Range("L107").Value = 696.5454599999999
$ws.Range("N107").Value = -4536.54546
$ws.Range("H113").Value = 3320
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 4093.3333
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 12279.9999
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -16619.9999
$ws.Range("H135").Value = 499.04443
$ws.Range("I135").Value = 426
$ws.Range("J135").Value = 836.875
$ws.Range("K135").Value = 3834
$ws.Range("L135").Value = 7531.875
$ws.Range("M135").Value = -1299
$ws.Range("N135").Value = -12601.875
$ws.Range("H140").Value = 1774.4445
$ws.Range("J140").Value = 1912.5
$ws.Range("L140").Value = 5737.5
$ws.Range("N140").Value = -16097.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9000
$ws.Range("J5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("N5").Value = -9224
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3500.5
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 4001
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 4001
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -4225
$ws.Range("H126").Value = 3500.5
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 4001
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 12003
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -16943

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 12440
$ws.Range("J43").Value = 12440
$ws.Range("L43").Value = 12440
$ws.Range("N43").Value = -12738
$ws.Range("H97").Value = 33500
$ws.Range("J97").Value = 33500
$ws.Range("L97").Value = 33500
$ws.Range("N97").Value = -35482
$ws.Range("H122").Value = 1430.5238
$ws.Range("I122").Value = 1127.625
$ws.Range("J122").Value = 2399.8
$ws.Range("K122").Value = 3382.875
$ws.Range("L122").Value = 7199.400000000001
$ws.Range("M122").Value = -932.875
$ws.Range("N122").Value = -12099.4
